# Applies the small textual corrections described by the commit:
#  - "Hep B Immunization (among 1 year olds (%))" -> "Hepatitis B Immunization (among 1-year-olds (%))"
#  - "Polio Immunization (among 1 year olds (%))"  -> "Polio Immunization (among 1-year-olds (%))"
#  - "Diphtheria Immunization (among 1 year olds (%))" -> "Diphtheria Immunization (among 1-year-olds (%))"
#  - "HIV/AIDS (cases per 1000)" -> "HIV/AIDS (deaths per 1000 live births)"
#  - "Measles (cases per 1000)"  -> "Measles (reported cases per 1000)"
#  - stray trailing "git" removed from the last bullet

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Hep B Immunization (among 1 year olds (%)) vs. Life Expectancy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Hepatitis B Immunization (among 1-year-olds (%)) vs. Life Expectancy",
    2)

$d.Content.Find.Execute(
    "Polio Immunization (among 1 year olds (%)) vs. Life Expectancy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Polio Immunization (among 1-year-olds (%)) vs. Life Expectancy",
    2)

$d.Content.Find.Execute(
    "Diphtheria Immunization (among 1 year olds (%)) vs. Life Expectancy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Diphtheria Immunization (among 1-year-olds (%)) vs. Life Expectancy",
    2)

$d.Content.Find.Execute(
    "HIV/AIDS (cases per 1000) vs. Life Expectancy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "HIV/AIDS (deaths per 1000 live births) vs. Life Expectancy",
    2)

$d.Content.Find.Execute(
    "Measles (cases per 1000) vs. Life Expectancy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Measles (reported cases per 1000) vs. Life Expectancy",
    2)

$d.Content.Find.Execute(
    "NaN = 0 r-value: - .14git",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NaN = 0 r-value: - .14",
    2)
